$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.179.07"
$ws.Range("E2").Value = "  +1.13%  "

$ws.Range("D3").Value = "1.642.05"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.20"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.515"
$ws.Range("E6").Value = "  +1.38%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.256"
$ws.Range("E8").Value = "  +0.94%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.03"
$ws.Range("E10").Value = "  +1.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").Value = "1.871.29"
$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").Value = "1.642.71"
$ws.Range("E13").Value = "  -0.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.16"
$ws.Range("E14").Value = "  +0.63%  "

$ws.Range("E15").Value = "  +2.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.35"
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("D17").Value = "27.150.92"
$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("E18").Value = "  +1.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.68"
$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.88"
$ws.Range("E21").Value = "  +3.70%  "

$ws.Range("E22").Value = "  +6.50%  "

$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.69"

$ws.Range("E26").Value = "  +1.86%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.119"
$ws.Range("E28").Value = "  -0.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.79"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0510"
$ws.Range("E30").Value = "  +0.36%  "

$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("E32").Value = "  +0.83%  "

$ws.Range("E33").Value = "  +1.01%  "

$ws.Range("E34").Value = "  +0.95%  "

$ws.Range("D35").Value = "1.272.89"
$ws.Range("E35").Value = "  +2.28%  "

$ws.Range("E37").Value = "  +1.87%  "

$ws.Range("E38").Value = "  +2.72%  "

$ws.Range("E39").Value = "  +0.86%  "

$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  +7.16%  "

$ws.Range("E43").Value = "  -1.34%  "

$ws.Range("D44").Value = "1.782.39"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.80"
$ws.Range("E45").Value = "  +1.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.89"
$ws.Range("E46").Value = "  +0.45%  "

$ws.Range("E47").Value = "  +1.70%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  -0.14%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0975"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").Value = "  +0.37%  "
